$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty rows with the new section entries.
$ws.Range("A28").Value = "4.1 Objects"
$ws.Range("A30").Value = "4.3 Symbols skipp for now use a reference"
$ws.Range("A34").Value = "5.1 Methods of primitives"

# Update the selected cell to match the committed state.
$ws.Range("A30").Select()
